$newCharacteristics = @(
    "Ammonium",
    "Chl a",
    "Chl a (probe)",
    "Chloride",
    "Chlorophyll a",
    "Chlorophyll a (probe)",
    "Chlorophyll a (probe) concentration, Cyanobacteria (bluegreen)",
    "Conductivity",
    "Cyanobacteria (lab)",
    "Cyanobacteria (probe)",
    "Depth",
    "Depth, Secchi disk depth",
    "Dissolved oxygen (DO)",
    "Dissolved oxygen saturation",
    "DO",
    "DO saturation",
    "E.coli",
    "Enterococcus",
    "Escherichia coli",
    "Fecal Coliform",
    "Flow",
    "Gage",
    "Height, gage",
    "Metals",
    "Microcystins",
    "Nitrate",
    "Nitrate + Nitrite",
    "Nitrite",
    "Ortho P",
    "Orthophosphate",
    "Particulate organic carbon",
    "pH",
    "Pheophytin",
    "Pheophytin a",
    "Phosphorus, Particulate Organic",
    "POC",
    "PON",
    "POP",
    "Salinity",
    "Secchi Depth",
    "Silicate",
    "Sp Conductance",
    "Specific conductance",
    "Sulfate",
    "Surfactants",
    "TDN",
    "TDP",
    "TDS",
    "Temperature, air",
    "Temperature, water",
    "TKN",
    "TN",
    "Total dissolved solids",
    "Total Kjeldahl nitrogen",
    "Total Nitrogen, mixed forms",
    "Total Phosphorus, mixed forms",
    "Total suspended solids",
    "TP",
    "TSS",
    "Turbidity",
    "Water Temp"
)

$wb = $excel.ActiveWorkbook

# --- Values sheet: update Characteristic Name pick list (column D) ---
$wsValues = $wb.Worksheets.Item("Values")
for ($i = 0; $i -lt $newCharacteristics.Length; $i++) {
    $row = 5 + $i
    $wsValues.Cells.Item($row, 4).Value = $newCharacteristics[$i]
}

# Remove now-unused rows 66:87 (old list had duplicates / extra trailing rows)
$wsValues.Rows("66:87").Delete()

# Re-apply sort so the sortState reference shrinks to the new range
$wsValues.Sort.SortFields.Clear()
$wsValues.Sort.SortFields.Add($wsValues.Range("D2:D65")) | Out-Null
$wsValues.Sort.SetRange($wsValues.Range("D2:D65"))
$wsValues.Sort.Header = 0
$wsValues.Sort.Apply()

$wsValues.Range("D65").Select()

# --- Results sheet: update data validation list range + selection ---
$wsResults = $wb.Worksheets.Item("Results")
$rngH = $wsResults.Range("H2:H1048576")
$rngH.Validation.Delete()
$rngH.Validation.Add(3, 1, 1, "=Values!`$D`$2:`$D`$65")
$rngH.Validation.InCellDropdown = $true

# --- Instructions sheet: update template date note ---
$wsInstructions = $wb.Worksheets.Item("Instructions")
$wsInstructions.Range("C1").Value = "Template updated 6/29/23"

# --- Restore selection / active sheet state ---
$wsResults.Activate()
$wsResults.Range("A6").Select()
